$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) columns for rows 2-11
$ws.Range("C2").Value = 0.7362970845300184
$ws.Range("D2").Value = 0.4665999738702449

$ws.Range("C3").Value = 0.8202600933124645
$ws.Range("D3").Value = 0.4177816105352847

$ws.Range("C4").Value = 0.2232295455291295
$ws.Range("D4").Value = 0.8246926429841099

$ws.Range("C5").Value = 0.9643086159758357
$ws.Range("D5").Value = 0.3417002327762781

$ws.Range("C6").Value = -0.009251586009557056
$ws.Range("D6").Value = 0.9926724803319931

$ws.Range("C7").Value = -0.7591752114187625
$ws.Range("D7").Value = 0.4529792443776581

$ws.Range("C8").Value = 0.1703466144765603
$ws.Range("D8").Value = 0.8657471940528998

$ws.Range("C9").Value = -0.616056509202926
$ws.Range("D9").Value = 0.5419599843019625

$ws.Range("C10").Value = 0.1387179134870179
$ws.Range("D10").Value = 0.8904914404532021

$ws.Range("C11").Value = 0.6700993964502581
$ws.Range("D11").Value = 0.5073205123181901
